# Rotates the data rows 2-9 of the "Artfynd" sheet.
#
# Observed effect (per row-id column A):
#   row 2 <- old row 9
#   row 3 <- old row 3   (unchanged)
#   row 4 <- old row 2
#   row 5 <- old row 4
#   row 6 <- old row 5
#   row 7 <- old row 6
#   row 8 <- old row 7
#   row 9 <- old row 8
#
# i.e. the record that used to live in row 9 moves up to row 2, and every
# other record (rows 2,4,5,6,7,8) shifts down by one row, wrapping back to
# row 9 via row 8. Row 3 keeps its original content.
#
# Implemented generically: read the full contents of every used column for
# rows 2-9 into memory first (capturing the literal text for text columns so
# values such as "80" or "2021-08-04" are not misread as numbers/dates),
# then write the rotated data back so no record's data is overwritten
# before it has been captured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are genuinely numeric in this sheet.
$numericCols = @("A","B","E","Q","R","S")
# Columns whose values are booleans (TRUE/FALSE) in this sheet.
$boolCols = @("AD","AE","AG")
# Everything else is free text (some of it, e.g. "80" or "2021-08-04",
# would otherwise be auto-detected as a number/date by Excel).
$textCols = @("C","D","F","G","H","I","J","K","L","M","N","P","T","U","V","W","Y","Z","AA","AB","AT","AW","AX","AY")

$allCols = $numericCols + $boolCols + $textCols
$dataRows = @(2,3,4,5,6,7,8,9)

function Get-CellSnapshot($rng) {
    $v = $rng.Value()
    if ($v -eq $null) { return $null }
    # Use the literal displayed text so numeric-looking / date-looking
    # strings round-trip as text instead of being reinterpreted.
    return $rng.Text
}

function Set-TextCell($rng, $val) {
    if ($val -eq $null) {
        $rng.Value = $null
        return
    }
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# 1) Snapshot every cell in the affected rows.
$snapshot = @{}
foreach ($r in $dataRows) {
    $rowValues = @{}
    foreach ($col in $numericCols) {
        $rowValues[$col] = $ws.Range("$col$r").Value()
    }
    foreach ($col in $boolCols) {
        $rowValues[$col] = $ws.Range("$col$r").Value()
    }
    foreach ($col in $textCols) {
        $rowValues[$col] = Get-CellSnapshot($ws.Range("$col$r"))
    }
    $snapshot[$r] = $rowValues
}

# 2) Destination row -> source row mapping (the rotation described above).
$sourceForDest = @{2=9; 3=3; 4=2; 5=4; 6=5; 7=6; 8=7; 9=8}

# 3) Write the rotated data back out.
foreach ($destRow in $dataRows) {
    $srcRow = $sourceForDest[$destRow]
    $srcValues = $snapshot[$srcRow]

    foreach ($col in $numericCols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
    foreach ($col in $boolCols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
    foreach ($col in $textCols) {
        Set-TextCell $ws.Range("$col$destRow") $srcValues[$col]
    }
}
